$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 changes
$ws.Range("V2").Value = 1.5
$ws.Range("W2").Value = 1.48

# Row 3 changes
$ws.Range("F3").Value = 1.92
$ws.Range("G3").Value = 2.62
$ws.Range("H3").Value = 1.42
$ws.Range("J3").Value = 2.48
$ws.Range("K3").Value = 4.2
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 1.33
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.16
$ws.Range("S3").Value = 1.01
$ws.Range("W3").Value = 1.61

# Row 4 changes
$ws.Range("F4").Value = 2.32
$ws.Range("I4").Value = 3.4
$ws.Range("K4").Value = 4.1
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 1.7
